$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=293; A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=44628; E=13; F="Fruta"; G=100107; H="Otros"; I=100107011; J="Tuna"; K="Sin especificar"; L="Especial";  M=200; N=15000; O=15000; P=15000; Q="$/caja 18 kilos"; R="Región Metropolitana"; S=833; T=18 },
    @{ Row=294; A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=44628; E=13; F="Fruta"; G=100107; H="Otros"; I=100107011; J="Tuna"; K="Sin especificar"; L="Primera";   M=200; N=12000; O=12000; P=12000; Q="$/caja 18 kilos"; R="Región Metropolitana"; S=667; T=18 },
    @{ Row=295; A=6; B="Mercado Mayorista Lo Valledor de Santiago"; C="Metropolitana"; D=44628; E=13; F="Fruta"; G=100107; H="Otros"; I=100107011; J="Tuna"; K="Sin especificar"; L="Segunda";   M=200; N=10000; O=10000; P=10000; Q="$/caja 18 kilos"; R="Región Metropolitana"; S=556; T=18 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($n, 5).Value = $r.E
    $ws.Cells.Item($n, 6).Value = $r.F
    $ws.Cells.Item($n, 7).Value = $r.G
    $ws.Cells.Item($n, 8).Value = $r.H
    $ws.Cells.Item($n, 9).Value = $r.I
    $ws.Cells.Item($n, 10).Value = $r.J
    $ws.Cells.Item($n, 11).Value = $r.K
    $ws.Cells.Item($n, 12).Value = $r.L
    $ws.Cells.Item($n, 13).Value = $r.M
    $ws.Cells.Item($n, 14).Value = $r.N
    $ws.Cells.Item($n, 15).Value = $r.O
    $ws.Cells.Item($n, 16).Value = $r.P
    $ws.Cells.Item($n, 17).Value = $r.Q
    $ws.Cells.Item($n, 18).Value = $r.R
    $ws.Cells.Item($n, 19).Value = $r.S
    $ws.Cells.Item($n, 20).Value = $r.T
}
